# ============================================================================
# Applies the bdstricture01 ERCP-report templating edit.
# ============================================================================
$d = $word.ActiveDocument
$VT = [char]11   # vertical-tab -> becomes a <w:br/> inside a run when assigned to Range.Text

# ----------------------------------------------------------------------------
# 1. Indications paragraph: rewrite the body text.
# ----------------------------------------------------------------------------
$d.Paragraphs(3).Range.Text = "70-year-old male patient is here for an ERCP procedure for management of a malignant distal biliary stricture."

# Insert the new "Medications" / "Monitoring" subsections right after it.
$r = $d.Paragraphs(3).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(4)
$np.Style = "Heading3"
$np.Range.Text = "Medications"

$r = $np.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(5)
$np.Style = "Normal"
$np.Range.Text = "Refer to record of source."

$r = $np.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(6)
$np.Style = "Heading3"
$np.Range.Text = "Monitoring"

$r = $np.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(7)
$np.Style = "Normal"
$np.Range.Text = "Johns Hopkins Standard."

# ----------------------------------------------------------------------------
# 2. "EGD Findings" -> "History" heading, body text rewritten.
#    (was paragraphs 4/5, now shifted down by 4 -> 8/9)
#    Paragraph 9 is a multi-run paragraph (bold "ESOPHAGUS:"/"STOMACH:" runs),
#    so a plain Range.Text assignment would only touch the first run and
#    leave the rest behind. Instead insert a brand-new paragraph with the
#    desired text immediately before it, then delete the old one outright.
# ----------------------------------------------------------------------------
$d.Paragraphs(8).Range.Text = "History"

$historyText = "The patient reports  history of management of a malignant distal biliary stricture." + $VT + `
    "Patient has a history of recent diagnosis of metastatic pancreatic cancer discovered during a routine AAA evaluation. Mass is located in the head of the pancreas and it is associated with peripancreatic and retroperitoneal lymphadenopathy with lung nodules." + $VT + `
    "Recent diagnosis of metastatic pancreatic cancer with mass in the head of the pancreas, peripancreatic and retroperitoneal lymphadenopathy, lung nodules. Current medications include undergoing chemotherapy as part of a clinical trial."

$r = $d.Paragraphs(8).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(9)
$np.Style = "Normal"
$np.Range.Text = $historyText
# Remove the stale multi-run EGD paragraph, now pushed down to index 10.
$d.Paragraphs(10).Range.Delete()

# ----------------------------------------------------------------------------
# 3. "ERCP Findings" -> "Description of Procedure" heading, body rewritten.
#    (was paragraphs 6/7, now shifted -> 10/11)
# ----------------------------------------------------------------------------
$d.Paragraphs(10).Range.Text = "Description of Procedure"
$d.Paragraphs(11).Range.Text = "After the risks, benefits and alternatives of the procedure were thoroughly explained, informed consent was verified, confirmed and timeout was successfully executed by the treatment team. With the patient in the left semi-prone position, medications were administered intravenously. The duodenoscope Olympus TJF Q180V was passed from the mouth into the esophagus and further advanced from the esophagus into the stomach. From the stomach, the scope was directed to the second portion of the duodenum."

# ----------------------------------------------------------------------------
# 4. Insert a new "Findings" section right after Description of Procedure.
# ----------------------------------------------------------------------------
$r = $d.Paragraphs(11).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(12)
$np.Style = "Heading2"
$np.Range.Text = "Findings"

$r = $np.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(13)
$np.Style = "Normal"
$findingsText = "A scout film of the abdomen was performed.  plastic biliary stent in situ. A scout film of the abdomen reveals a plastic biliary stent in situ.." + $VT + $VT + `
    "The duodenoscope Olympus TJF Q180V was advanced to the second portion of the duodenum without difficulty and without detailed examination of the upper GI tract." + $VT + `
    " The esophagus, stomach, and duodenum appeared unremarkable on limited inspection." + $VT + `
    "The major papilla was identified. It appeared normal in position and morphology, with an intact orifice and no surrounding erythema or edema.A previously placed plastic or metal stent was seen emerging from the major papilla." + $VT + `
    "The minor papilla could not be located." + $VT + $VT + `
    "The ampulla was identified and found to be normal." + $VT + $VT + `
    "Bile duct cannulation was attempted using a sphincterotome preloaded with a guidewire (0.035 inch guide wire)." + $VT + `
    "Bile duct cannulation was successful. Biliary cannulation was achieved without pancreatic duct entry." + $VT + $VT + `
    "Contrast was injected under fluoroscopic guidance and cholangiogram was performed." + $VT + `
    "The common bile duct (CBD) measured 14.0 mm. Intrahepatic ducts were mildly dilated. 3 centimeter biliary stricture in the distal part. Occluded plastic biliary stent is identified and removed using a snare. Final cholangiogram confirms persistent obstruction of the proximal mid-CBD.." + $VT + `
    "Pancreatogram was not performed." + $VT + $VT + `
    "Sphincterotomy was performed." + $VT + $VT + `
    "A 10-60 mm uncovered self-expandable metal stent deployed across the distal biliary stricture was placed. The scope was then completely withdrawn from the patient and the procedure completed." + $VT + $VT + `
    "Estimated blood loss: None." + $VT + $VT + `
    "Specimens removed: None." + $VT + $VT + `
    "Complications: There were no immediate complications."
$np.Range.Text = $findingsText

# ----------------------------------------------------------------------------
# 5. Insert a new "ERCP Quality Metrics" section right after Findings.
# ----------------------------------------------------------------------------
$r = $d.Paragraphs(13).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(14)
$np.Style = "Heading2"
$np.Range.Text = "ERCP Quality Metrics"

$r = $np.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$np = $d.Paragraphs(15)
$np.Style = "Normal"
$metricsText = "Difficulty of ERCP:  pancreatic duct cannulation not attempted." + $VT + `
    "Cannulation success: Yes." + $VT + `
    "Post-ERCP pancreatitis prophylaxis:No." + $VT + `
    "Rectal indomethacin: No." + $VT + `
    "Successful completion of intended procedure: Yes." + $VT + `
    "Failed ERCP from another facility or provider: No."
$np.Range.Text = $metricsText

# ----------------------------------------------------------------------------
# 6. Impressions list: rewrite items 1-4, delete old items 5 & 6.
#    (Impressions heading now at 16, items 17..22)
# ----------------------------------------------------------------------------
$d.Paragraphs(17).Range.Text = "1. Successful ERCP with stent exchange and sphincterotomy"
$d.Paragraphs(18).Range.Text = "2. 3 cm distal biliary stricture with upstream common bile duct dilation up to 14 mm"
$d.Paragraphs(19).Range.Text = "3. Mild diffuse intrahepatic ductal dilation"
$d.Paragraphs(20).Range.Text = "4. 10-60 mm uncovered self-expandable metal stent deployed across the distal biliary stricture"
# Remove old items 5 and 6 (now paragraphs 21 and 21 again after first delete)
$d.Paragraphs(21).Range.Delete()
$d.Paragraphs(21).Range.Delete()

# ----------------------------------------------------------------------------
# 7. Recommendations list: item 3 replaced by the old item 4's text, drop the
#    trailing 4th item.
#    (Recommendations heading at 21, items 22..25)
# ----------------------------------------------------------------------------
$d.Paragraphs(24).Range.Text = "3. Follow up with referring provider."
$d.Paragraphs(25).Range.Delete()
